# holly added S.GISH to harvester in bioSamples -- fixes the harvester column
# in rnaSample_0673: column B ("harvester") previously held "Retrofitted_0673"
# for every data row; it should now read "S.GISH".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the harvester column (B2:B16) for all 15 data rows.
$ws.Range("B2:B16").Value = "S.GISH"

# Widen column B slightly (matches the author's manual column resize after
# editing, raw OOXML width ~8.83; ColumnWidth property uses a slightly
# different unit, so back it off by the standard padding offset).
$ws.Columns("B").ColumnWidth = 7.96

# Reflect that column B was the focus of the edit (active cell B1, whole
# column selected), matching the saved selection state.
$null = $ws.Range("B:B").Select()
